# Anonymise the "Early Detector" results sheet: replace real GitHub
# usernames in column A with generic "UserN" placeholders. Each of the
# three project teams (P1, P2, P3) gets its own User1..User5 numbering
# based on the member's position in that team's roster (P1's 4th and
# 5th members both collapse onto "User4", matching the original data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# P1 team members
$ws.Cells.Replace("SohanK", "User1")
$ws.Cells.Replace("anishashetty", "User2")
$ws.Cells.Replace("rroycho", "User3")
$ws.Cells.Replace("anlawande", "User4")
$ws.Cells.Replace("ptrived", "User4")

# P2 team members
$ws.Cells.Replace("krishnatejadinavahi", "User1")
$ws.Cells.Replace("keemen90", "User2")
$ws.Cells.Replace("smruthiEJ", "User3")
$ws.Cells.Replace("RonakNisher", "User4")
$ws.Cells.Replace("juhidesai", "User5")

# P3 team members
$ws.Cells.Replace("bhashwanth", "User1")
$ws.Cells.Replace("kumar-utsav", "User2")
$ws.Cells.Replace("rarora4", "User3")
$ws.Cells.Replace("yatish27", "User4")

# Restore the author's last-saved selection on the sheet.
$ws.Range("F5").Select()
